$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.116.48'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '1.795.75'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'308.88"
$ws.Range('E5').Value = '  -2.02%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = "'1.008"
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').Value = "'0.4180"
$ws.Range('E7').Value = '  -1.50%  '
$ws.Range('D8').Value = "'0.3556"
$ws.Range('E8').Value = '  -3.62%  '
$ws.Range('D9').Value = "'0.07076"
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('D10').Value = "'0.8439"
$ws.Range('E10').Value = '  -2.44%  '
$ws.Range('D11').Value = "'20.07"
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('D12').Value = '1.835.96'
$ws.Range('E12').Value = '  -13.47%  '
$ws.Range('D13').Value = "'5.280"
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'6.335"
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').Value = "'0.06856"
$ws.Range('E15').Value = '  -2.02%  '
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = "'79.70"
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = "'0.000008719"
$ws.Range('E18').Value = '  -3.07%  '
$ws.Range('D19').Value = "'1.009"
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').Value = "'15.05"
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('D21').Value = '27.394.63'
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').Value = "'5.052"
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').Value = "'10.74"
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').Value = '2.070.01'
$ws.Range('E24').Value = '  -2.99%  '
$ws.Range('D25').Value = "'1.965"
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('D26').Value = "'153.36"
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').Value = "'18.14"
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').Value = "'5.022"
$ws.Range('E28').Value = '  -4.24%  '
$ws.Range('D29').Value = "'112.58"
$ws.Range('E29').Value = '  -3.02%  '
$ws.Range('D30').Value = "'1.661"
$ws.Range('E30').Value = '  -9.57%  '
$ws.Range('D31').Value = "'0.08896"
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'0.7249"
$ws.Range('E32').Value = '  -5.72%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = "'2.894"
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').Value = "'4.366"
$ws.Range('E34').Value = '  -4.06%  '
$ws.Range('D35').Value = "'1.007"
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').Value = "'1.074"
$ws.Range('E36').Value = '  -5.56%  '
$ws.Range('D37').Value = "'1.071"
$ws.Range('E37').Value = '  -3.28%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.05108"
$ws.Range('E38').Value = '  -4.56%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.01892"
$ws.Range('E39').Value = '  -3.40%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = "'0.1620"
$ws.Range('E40').Value = '  -2.38%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.4944"
$ws.Range('E41').Value = '  -2.87%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = "'2.666"
$ws.Range('E42').Value = '  -5.66%  '
$ws.Range('D43').Value = "'6.244"
$ws.Range('E43').Value = '  -8.78%  '
$ws.Range('D44').Value = "'8.031"
$ws.Range('E44').Value = '  -5.30%  '
$ws.Range('D45').Value = "'1.008"
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').Value = "'104.65"
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('D48').Value = "'0.06312"
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('D49').Value = "'0.4532"
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').Value = "'1.590"
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('D51').Value = "'62.16"
$ws.Range('E51').Value = '  -2.82%  '
